# NIT-9011175755.xlsx — "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
#
# Logic (verified by diffing the decoded before/after tables):
#   - Rows 16-18 (period 2507 entries for EDWIN DELAS AGUAS PERTUZ,
#     JOSE ANDRES CASTILLA ROMERO and DUBAN ENRIQUE ESTRADA MENDOZA) are
#     removed outright (old EC rows).
#   - The remaining period-2507 row (ARLINSON PEREZ RECUERO) shifts up to
#     become the new row 16, the two period-2508 rows (JOSE / DUBAN) shift
#     up to rows 17-18, and the trailing period-2508 row for ARLINSON PEREZ
#     RECUERO shifts up to row 19 and has its period bumped from 2508 to
#     2509 (new EC row).
#   - The footer block (signature lines) that used to sit at rows 27-28
#     naturally re-flows to rows 24-25 once the 3 rows are deleted.
#   - The summary cells (VALOR MORA total, worker count, period count) are
#     refreshed to match the new 4-row table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Delete the three obsolete rows (old period-2507 entries for EDWIN,
#    JOSE and DUBAN). This shifts every row below up by 3, which is exactly
#    what turns the old B2:J28 sheet into the new B2:J25 sheet and moves the
#    footer from rows 27/28 to rows 24/25 automatically.
$ws.Rows("16:18").Delete()

# 2) The row that used to be "ARLINSON PEREZ RECUERO / 2508" (old row 22)
#    is now row 19 - bump its period from 2508 to 2509 (new EC entry).
$ws.Range("E19").Value = "2509"

# 3) Refresh the summary header cells to match the new 4-row table.
#    VALOR MORA total = sum of the Valor Mora column (F16:F19).
$ws.Range("E11").Value = 430640
#    Cant. Trabajadores (distinct workers) and Cant. Periodos (distinct
#    periods) for the refreshed table.
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 3
